$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 155.157341903
$ws.Range("C2").Value = 0.19759461697799999

$ws.Range("B3").Value = 155157.34190299999
$ws.Range("C3").Value = 5927.8385093399993

$ws.Range("B4").Value = 573104.67378711107
$ws.Range("C4").Value = 27240.492693895569

$ws.Range("B5").Value = 11462.09347574222
$ws.Range("C5").Value = 544.80985387791134
